$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.599.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.587.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.62%  "

$ws.Range("E6").Value = "  -3.22%  "

$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.593.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.64%  "

$ws.Range("E11").Value = "  -2.30%  "

$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("E13").Value = "  +1.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.040.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.547.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.64%  "

$ws.Range("E17").Value = "  -0.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.589.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.51%  "

$ws.Range("E19").Value = "  -1.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.17%  "

$ws.Range("E22").Value = "  -1.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.56%  "

$ws.Range("E25").Value = "  -1.65%  "

$ws.Range("E26").Value = "  -1.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.702.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.18%  "

$ws.Range("E29").Value = "  -2.70%  "

$ws.Range("E30").Value = "  -3.00%  "

$ws.Range("E31").Value = "  +0.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "152.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.40%  "

$ws.Range("E34").Value = "  -2.09%  "

$ws.Range("E35").Value = "  +1.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("E37").Value = "  -2.90%  "

$ws.Range("E38").Value = "  +4.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.848"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.30%  "

$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "296.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.623"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0995"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0559"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.61%  "

$ws.Range("E47").Value = "  +0.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.17%  "

$ws.Range("E50").Value = "  -2.72%  "

$ws.Range("E51").Value = "  +0.14%  "
